$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the wage & contribution column headers for consistency.
# (A1 "year", B1 "anticipated wages" and I1 "big-ticket items" keep their text.)
$ws.Range("H1").Value = "Roth conv"
$ws.Range("C1").Value = "taxable ctrb"
$ws.Range("D1").Value = "401k ctrb"
$ws.Range("E1").Value = "Roth 401k ctrb"
$ws.Range("F1").Value = "IRA ctrb"
$ws.Range("G1").Value = "Roth IRA ctrb"

# Bold the header row.
$ws.Range("A1:I1").Font.Bold = $true

# Select the header row (row 1) instead of row 2.
$ws.Rows("1:1").Select()
